# MoviesList.xlsx update:
#  - replace the two movie rows (url + name) with new entries
#  - turn the url cells into real hyperlinks (adds the built-in "Hyperlink"
#    style/font, matching Excel's own "Insert Hyperlink" behaviour)
#  - restore page setup (paper size / orientation) and move the selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Táo Quậy ---------------------------------------------------
$ws.Range("A2").Value = "https://1080.opstream4.com/20220309/1217_883067c1/1000k/hls/mixed.m3u8"
$ws.Range("B2").Value = "Táo Quậy"

# --- Row 3: Thang Máy ---------------------------------------------------
$ws.Range("A3").Value = "https://hd.hdbophim.com/20221116/26035_8ae84283/1163617/hls/mixed.m3u8"
$ws.Range("B3").Value = "Thang Máy"

# Turn the two URL cells into clickable hyperlinks
$ws.Hyperlinks.Add($ws.Range("A2"), $ws.Range("A2").Value2)
$ws.Hyperlinks.Add($ws.Range("A3"), $ws.Range("A3").Value2)

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move the active selection to A4, as recorded in the saved view state
$ws.Range("A4").Select()
